# narrative for fig 3 - Findings page
# Updates the lower/upper confidence bound values (columns F and G) for
# several rows on the Findings page data table, per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 0.78

$ws.Range("F14").Value = 0.58796992481203
$ws.Range("G14").Value = 0.7082706766917293

$ws.Range("F15").Value = 0.5714285714285714
$ws.Range("G15").Value = 0.6932706766917288

$ws.Range("F16").Value = 0.5714285714285714

$ws.Range("F17").Value = 0.5714285714285714

$ws.Range("F18").Value = 0.8322222222222222
$ws.Range("G18").Value = 0.9066666666666666

$ws.Range("G19").Value = 0.8833333333333333

$ws.Range("F20").Value = 0.7988888888888889
$ws.Range("G20").Value = 0.8788888888888888

$ws.Range("G21").Value = 0.8744444444444445

$ws.Range("F22").Value = 0.8913043478260869

$ws.Range("G23").Value = 0.9608695652173913

$ws.Range("F24").Value = 0.8717391304347826

$ws.Range("F25").Value = 0.8695108695652174
